$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "police"
$ws.Range("A3").Value = "schools"
$ws.Range("C3").Value = 10
